$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$val0 = @'
Unknown Title
'@
$ws.Range("C2").Value = $val0

$val1 = @'
Unknown Abstract
'@
$ws.Range("D2").Value = $val1

$val2 = @'
[]
'@
$ws.Range("E2").Value = $val2

$val3 = @'
not found
'@
$ws.Range("F2").Value = $val3

$val4 = @'
N/A
'@
$ws.Range("G2").Value = $val4

$val5 = @'
'1970-01-01
'@
$ws.Range("H2").Value = $val5

$val6 = @'
Unknown Title
'@
$ws.Range("C3").Value = $val6

$val7 = @'
Unknown Abstract
'@
$ws.Range("D3").Value = $val7

$val8 = @'
[]
'@
$ws.Range("E3").Value = $val8

$val9 = @'
not found
'@
$ws.Range("F3").Value = $val9

$val10 = @'
N/A
'@
$ws.Range("G3").Value = $val10

$val11 = @'
'1970-01-01
'@
$ws.Range("H3").Value = $val11

$val12 = @'
Unknown Title
'@
$ws.Range("C4").Value = $val12

$ws.Range("J4").Value = ""
$val14 = @'
[Chaolin%Huang%NULL%0, Yeming%Wang%NULL%0, Xingwang%Li%NULL%0, Lili%Ren%NULL%0, Jianping%Zhao%NULL%0, Yi%Hu%NULL%0, Li%Zhang%NULL%0, Guohui%Fan%NULL%0, Jiuyang%Xu%NULL%0, Xiaoying%Gu%NULL%0, Zhenshun%Cheng%NULL%0, Ting%Yu%NULL%0, Jiaan%Xia%NULL%0, Yuan%Wei%NULL%0, Wenjuan%Wu%NULL%0, Xuelei%Xie%NULL%0, Wen%Yin%NULL%0, Hui%Li%NULL%0, Min%Liu%NULL%0, Yan%Xiao%NULL%0, Hong%Gao%NULL%0, Li%Guo%NULL%0, Jungang%Xie%NULL%0, Guangfa%Wang%NULL%0, Rongmeng%Jiang%NULL%0, Zhancheng%Gao%NULL%0, Qi%Jin%NULL%0, Jianwei%Wang%wangjw28@163.com%0, Bin%Cao%caobin_ben@163.com%0]
'@
$ws.Range("E5").Value = $val14

$ws.Range("I5").Value = ""
$val16 = @'
Elsevier Ltd.
'@
$ws.Range("J5").Value = $val16

$val17 = @'
[Yan%Deng%NULL%0, Wei%Liu%NULL%0, Kui%Liu%NULL%0, Yuan-Yuan%Fang%NULL%0, Jin%Shang%NULL%5, Ling%Zhou%NULL%0, Ke%Wang%NULL%0, Fan%Leng%NULL%5, Shuang%Wei%NULL%0, Lei%Chen%NULL%5, Hui-Guo%Liu%NULL%0, Pei-Fang%Wei%NULL%0, Pei-Fang%Wei%NULL%0]
'@
$ws.Range("E6").Value = $val17

$ws.Range("I6").Value = ""
$val19 = @'
Wolters Kluwer Health
'@
$ws.Range("J6").Value = $val19

$val20 = @'
In December 2019, a coronavirus 2019 (COVID-19) disease outbreak occurred in Wuhan, Hubei Province, China, and rapidly spread to other areas worldwide.
 Although diffuse alveolar damage and acute respiratory failure were the main features, the involvement of other organs needs to be explored.
 Since information on kidney disease in patients with COVID-19 is limited, we determined the prevalence of acute kidney injury (AKI) in patients with COVID-19. Further, we evaluated the association between markers of abnormal kidney function and death in patients with COVID-19. This was a prospective cohort study of 701 patients with COVID-19 admitted in a tertiary teaching hospital that also encompassed three affiliates following this major outbreak in Wuhan in 2020 of whom 113 (16.1%) died in hospital.
 Median age of the patients was 63 years (interquartile range, 50-71), including 367 men and 334 women.
 On admission, 43.9% of patients had proteinuria and 26.7% had hematuria.
 The prevalence of elevated serum creatinine, elevated blood urea nitrogen and estimated glomerular filtration under 60 ml/min/1.73m2 were 14.4, 13.1 and 13.1%, respectively.
 During the study period, AKI occurred in 5.1% patients.
 Kaplan-Meier analysis demonstrated that patients with kidney disease had a significantly higher risk for in-hospital death.
 Cox proportional hazard regression confirmed that elevated baseline serum creatinine (hazard ratio: 2.10, 95% confidence interval: 1.36-3.26), elevated baseline blood urea nitrogen (3.97, 2.57-6.14), AKI stage 1 (1.90, 0.76-4.76), stage 2 (3.51, 1.49-8.26), stage 3 (4.38, 2.31-8.31), proteinuria 1+ (1.80, 0.81-4.00), 2+∼3+ (4.84, 2.00-11.70), and hematuria 1+ (2.99, 1.39-6.42), 2+∼3+ (5.56,2.58- 12.01) were independent risk factors for in-hospital death after adjusting for age, sex, disease severity, comorbidity and leukocyte count.
 Thus, our findings show the prevalence of kidney disease on admission and the development of AKI during hospitalization in patients with COVID-19 is high and is associated with in-hospital mortality.
 Hence, clinicians should increase their awareness of kidney disease in patients with severe COVID-19.
'@
$ws.Range("D7").Value = $val20

$val21 = @'
[Yichun%Cheng%NULL%0, Ran%Luo%NULL%1, Kun%Wang%NULL%0, Meng%Zhang%NULL%3, Zhixiang%Wang%NULL%1, Lei%Dong%NULL%1, Junhua%Li%NULL%3, Ying%Yao%NULL%1, Shuwang%Ge%geshuwang@tjh.tjmu.edu.cn%1, Gang%Xu%xugang@tjh.tjmu.edu.cn%1]
'@
$ws.Range("E7").Value = $val21

$ws.Range("I7").Value = ""
$val23 = @'
International Society of Nephrology. Published by Elsevier Inc.
'@
$ws.Range("J7").Value = $val23

$val24 = @'
[Carly%Eastin%NULL%0, Travis%Eastin%NULL%2]
'@
$ws.Range("E8").Value = $val24

$ws.Range("I8").Value = ""
$val26 = @'
Elsevier
'@
$ws.Range("J8").Value = $val26

$val27 = @'
Unknown Title
'@
$ws.Range("C9").Value = $val27

$val28 = @'
[]
'@
$ws.Range("E9").Value = $val28

$val29 = @'
not found
'@
$ws.Range("F9").Value = $val29

$val30 = @'
N/A
'@
$ws.Range("G9").Value = $val30

$ws.Range("I9").Value = ""
$val32 = @'
[Manoocher%Soleimani%NULL%0]
'@
$ws.Range("E10").Value = $val32

$ws.Range("I10").Value = ""
$val34 = @'
MDPI
'@
$ws.Range("J10").Value = $val34

$val35 = @'
[Luwen%Wang%NULL%0, Xun%Li%NULL%3, Hui%Chen%NULL%0, Shaonan%Yan%NULL%2, Dong%Li%NULL%2, Yan%Li%NULL%0, Zuojiong%Gong%NULL%2]
'@
$ws.Range("E11").Value = $val35

$ws.Range("I11").Value = ""
$val37 = @'
S. Karger AG
'@
$ws.Range("J11").Value = $val37

$val38 = @'
[Xiao-Wei%Xu%NULL%0, Xiao-Xin%Wu%NULL%0, Xian-Gao%Jiang%NULL%0, Kai-Jin%Xu%NULL%0, Ling-Jun%Ying%NULL%0, Chun-Lian%Ma%NULL%0, Shi-Bo%Li%NULL%0, Hua-Ying%Wang%NULL%0, Sheng%Zhang%NULL%0, Hai-Nv%Gao%NULL%0, Ji-Fang%Sheng%NULL%0, Hong-Liu%Cai%NULL%0, Yun-Qing%Qiu%NULL%0, Lan-Juan%Li%NULL%0]
'@
$ws.Range("E12").Value = $val38

$ws.Range("I12").Value = ""
$val40 = @'
BMJ Publishing Group Ltd.
'@
$ws.Range("J12").Value = $val40

$val41 = @'
Unknown Title
'@
$ws.Range("C13").Value = $val41

$val42 = @'
Unknown Abstract
'@
$ws.Range("D13").Value = $val42

$val43 = @'
[]
'@
$ws.Range("E13").Value = $val43

$val44 = @'
not found
'@
$ws.Range("F13").Value = $val44

$val45 = @'
N/A
'@
$ws.Range("G13").Value = $val45

$val46 = @'
'1970-01-01
'@
$ws.Range("H13").Value = $val46

$val47 = @'
Background
id="Par1">Since December 2019, 2019 novel coronavirus pneumonia emerged in Wuhan city and rapidly spread throughout China and even the world.

 We sought to analyse the clinical characteristics and laboratory findings of some cases with 2019 novel coronavirus pneumonia .


Methods
id="Par2">In this retrospective study, we extracted the data on 95 patients with laboratory-confirmed 2019 novel coronavirus pneumonia in Wuhan Xinzhou District People’s Hospital from January 16th to February 25th, 2020. Cases were confirmed by real-time RT-PCR and abnormal radiologic findings.

 Outcomes were followed up until March 2th, 2020.
Results
id="Par3">Higher temperature, blood leukocyte count, neutrophil count, neutrophil percentage, C-reactive protein level, D-dimer level, alanine aminotransferase activity, aspartate aminotransferase activity, α - hydroxybutyrate dehydrogenase activity, lactate dehydrogenase activity and creatine kinase activity were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and so were lower lymphocyte count, lymphocyte percentage and total protein level.

 Age below 40 or above 60 years old, male, higher Creatinine level, and lower platelet count also seemed related to severe 2019 novel coronavirus pneumonia and composite endpoint, however the P values were greater than 0.05, which mean under the same condition studies of larger samples are needed in the future.


Conclusion
id="Par4">Multiple factors were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and more related studies are needed in the future.



'@
$ws.Range("D14").Value = $val47

$val48 = @'
[Gemin%Zhang%NULL%0, Jie%Zhang%945128911@qq.com%0, Bowen%Wang%NULL%3, Xionglin%Zhu%NULL%3, Qiang%Wang%NULL%6, Shiming%Qiu%NULL%3]
'@
$ws.Range("E14").Value = $val48

$ws.Range("I14").Value = ""
$val50 = @'
BioMed Central
'@
$ws.Range("J14").Value = $val50

$val51 = @'

              • COVID-19 has a great threat to world health.

'@
$ws.Range("D15").Value = $val51

$val52 = @'
[Xiaoli%Zhang%NULL%0, Huan%Cai%NULL%3, Jianhua%Hu%NULL%3, Jiangshan%Lian%NULL%3, Jueqing%Gu%NULL%3, Shanyan%Zhang%NULL%3, Chanyuan%Ye%NULL%0, Yingfeng%Lu%NULL%3, Ciliang%Jin%NULL%3, Guodong%Yu%NULL%3, Hongyu%Jia%NULL%3, Yimin%Zhang%NULL%3, Jifang%Sheng%jifang_sheng@zju.edu.cn%0, Lanjuan%Li%ljli@zju.edu.cn%0, Yida%Yang%yidayang65@zju.edu.cn%3]
'@
$ws.Range("E15").Value = $val52

$ws.Range("I15").Value = ""
$val54 = @'
The Author(s). Published by Elsevier Ltd on behalf of International Society for Infectious Diseases.
'@
$ws.Range("J15").Value = $val54

$val55 = @'
[Fei%Zhou%NULL%0, Ting%Yu%NULL%0, Ronghui%Du%NULL%0, Guohui%Fan%NULL%0, Ying%Liu%NULL%0, Zhibo%Liu%NULL%0, Jie%Xiang%NULL%0, Yeming%Wang%NULL%0, Bin%Song%NULL%0, Xiaoying%Gu%NULL%0, Lulu%Guan%NULL%0, Yuan%Wei%NULL%0, Hui%Li%NULL%0, Xudong%Wu%NULL%0, Jiuyang%Xu%NULL%0, Shengjin%Tu%NULL%0, Yi%Zhang%NULL%0, Hua%Chen%NULL%0, Bin%Cao%NULL%0]
'@
$ws.Range("E16").Value = $val55

$ws.Range("I16").Value = ""
$val57 = @'
Elsevier Ltd.
'@
$ws.Range("J16").Value = $val57

$val58 = @'
Unknown Title
'@
$ws.Range("C17").Value = $val58

$val59 = @'
Unknown Abstract
'@
$ws.Range("D17").Value = $val59

$val60 = @'
[]
'@
$ws.Range("E17").Value = $val60

$val61 = @'
not found
'@
$ws.Range("F17").Value = $val61

$val62 = @'
N/A
'@
$ws.Range("G17").Value = $val62

$val63 = @'
'1970-01-01
'@
$ws.Range("H17").Value = $val63

